$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.893.17"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "3.738.52"
$ws.Range("E3").Value = "  +6.91%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "420.43"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "131.82"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "3.730.60"
$ws.Range("E7").Value = "  +6.95%  "
$ws.Range("D8").Value = "0.647"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "0.772"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "0.184"
$ws.Range("E11").Value = "  +14.51%  "
$ws.Range("D12").Value = "0.0000408"
$ws.Range("E12").Value = "  +57.49%  "
$ws.Range("D13").Value = "42.95"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "10.55"
$ws.Range("E14").Value = "  +7.11%  "
$ws.Range("D15").Value = "4.313.85"
$ws.Range("E15").Value = "  +6.49%  "
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "20.83"
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("D18").Value = "3.722.93"
$ws.Range("E18").Value = "  +6.63%  "
$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("E20").Value = "  +4.84%  "
$ws.Range("D21").Value = "66.902.26"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Value = "445.69"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "16.46"
$ws.Range("D24").Value = "89.93"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "38.55"
$ws.Range("E26").Value = "  +13.95%  "
$ws.Range("D27").Value = "10.25"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").Value = "3.35"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "5.08"
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("D30").Value = "12.80"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +9.90%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "7.28"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "41.83"
$ws.Range("E35").Value = "  +5.81%  "
$ws.Range("D36").Value = "57.12"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "0.0496"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "0.0₃0758"
$ws.Range("E39").Value = "  +8.97%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.150"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  +30.04%  "
$ws.Range("D42").Value = "29.15"
$ws.Range("E42").Value = "  +35.11%  "
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "3.45"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").Value = "  +31.34%  "
$ws.Range("E46").Value = "  +6.54%  "
$ws.Range("D47").Value = "146.92"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "2.67"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").Value = "4.38"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("D51").Value = "0.308"
$ws.Range("E51").Value = "  -1.80%  "
